# Breakwater calculator update:
#  - Insert a new "Van Der Meer (1988a) ... slope 2.0:1" row into the
#    parameters table (computed from the existing 1.5:1 row via the
#    (2/1.5)^(1/3) cube scaling used elsewhere in this workbook), pushing
#    the table from 3 data rows to 4, and growing the "formulas"/"table"
#    named ranges to match.
#  - Add a 4th "Formula:" readout block (row 16) on the "formula" sheet
#    that mirrors rows 13-15, now that there is a 4th formula entry.

$wb = $excel.ActiveWorkbook
$paramWs  = $wb.Worksheets.Item("parameters")
$formWs   = $wb.Worksheets.Item("formula")

# --- 1. Grow the parameters table by one row -------------------------------
# Duplicate row 4's formatting (style ids + row height) down into row 5 so
# the new row looks identical to its neighbours, without touching every
# column on the sheet (a full-row Copy/Paste would stamp formatting across
# all 16384 columns).
$paramWs.Range("A4:F4").Copy()
$paramWs.Range("A5:F5").PasteSpecial(-4122)
$paramWs.Range("A5").EntireRow.RowHeight = $paramWs.Range("A4").EntireRow.RowHeight

# Shift the old row 3 (Chegini-Aghtouman, slope 1.5:1) down into row 5;
# row 4 (Chegini-Aghtouman, slope 2:1) is untouched.
$paramWs.Range("A5").Value = "Chegini-Aghtouman (2006) formula for Antifer Cubes, slope 1.5:1"
$paramWs.Range("B5").Value = 6.951
$paramWs.Range("C5").Value = 0.443
$paramWs.Range("D5").Value = 0.291
$paramWs.Range("E5").Value = 1.082
$paramWs.Range("F5").Value = 0.082

# Move the old row 2 (Van Der Meer, slope 1.5:1) down into row 3.
$paramWs.Range("A3").Value = "Van Der Meer (1988a) formula for Cubes (not Antifer), slope 1.5:1"
$paramWs.Range("B3").Value = 6.7
$paramWs.Range("C3").Value = 0.4
$paramWs.Range("D3").Value = 0.3
$paramWs.Range("E3").Value = 1
$paramWs.Range("F3").Value = 0.1

# Row 2 becomes the new Van Der Meer, slope 2.0:1 entry -- D50 and Hs/(D.D50)
# scale off the 1.5:1 numbers by the cube-root of the slope ratio, C/D/F stay.
$paramWs.Range("A2").Value = "Van Der Meer (1988a) formula for Cubes (not Antifer), slope 2.0:1"
$paramWs.Range("B2").Formula = "=6.7*(2/1.5)^(1/3)"
$paramWs.Range("C2").Value = 0.4
$paramWs.Range("D2").Value = 0.3
$paramWs.Range("E2").Formula = "=1*(2/1.5)^(1/3)"
$paramWs.Range("F2").Value = 0.1

# --- 2. Widen the named ranges that drive the "formula" sheet --------------
$wb.Names.Item("formulas").RefersTo = "=parameters!`$A`$2:`$A`$5"
$wb.Names.Item("table").RefersTo = "=parameters!`$B`$2:`$F`$5"

# --- 3. Add the 4th formula readout row on the "formula" sheet -------------
# Fill in the new row's content first, then stamp row 13's formatting (the
# "1=" / "Formula" / blank triple) on top so D16/E16/F16 end up styled like
# their siblings without the value write clobbering the pasted style.
$formWs.Range("D16").Value = "4="
$formWs.Range("E16").Formula = "=INDEX(formulas,4,1)"
$formWs.Range("D13:F13").Copy()
$formWs.Range("D16:F16").PasteSpecial(-4122)

$excel.CutCopyMode = 0
$wb.Application.Calculate()
